$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'c8daea83-fb7b-4236-b2d4-75329c45d9e1'
$ws.Range('B2').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D2').Value = 'Required column PAT_MRN_ID is missing in synthetic_fail.'
$ws.Range('H2').Value = 'Ensure synthetic_fail contains the column "PAT_MRN_ID"'
$ws.Range('A3').Value = '5231d332-934f-41a1-bf3a-bf4e729f91b7'
$ws.Range('B3').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D3').Value = 'Required column FACILITY is missing in synthetic_fail.'
$ws.Range('H3').Value = 'Ensure synthetic_fail contains the column "FACILITY"'
$ws.Range('A4').Value = '1f312a3c-27fb-432b-a263-e438e4185b71'
$ws.Range('B4').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D4').Value = 'Required column FIRST_NAME is missing in synthetic_fail.'
$ws.Range('H4').Value = 'Ensure synthetic_fail contains the column "FIRST_NAME"'
$ws.Range('A5').Value = '455b2050-381b-4964-8752-4143755e6996'
$ws.Range('B5').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D5').Value = 'Required column LAST_NAME is missing in synthetic_fail.'
$ws.Range('H5').Value = 'Ensure synthetic_fail contains the column "LAST_NAME"'
$ws.Range('A6').Value = 'f56834d4-40cf-478c-a4ba-bec92cc59265'
$ws.Range('B6').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D6').Value = 'Required column PAT_BIRTH_DATE is missing in synthetic_fail.'
$ws.Range('H6').Value = 'Ensure synthetic_fail contains the column "PAT_BIRTH_DATE"'
$ws.Range('A7').Value = '15868a63-eb31-47a5-9f55-0e71d000c03c'
$ws.Range('B7').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D7').Value = 'Required column MEDICAID_CIN is missing in synthetic_fail.'
$ws.Range('H7').Value = 'Ensure synthetic_fail contains the column "MEDICAID_CIN"'
$ws.Range('A8').Value = 'b9d0988e-afd8-4443-9df7-0dd02884df1a'
$ws.Range('B8').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D8').Value = 'Required column ENCOUNTER_ID is missing in synthetic_fail.'
$ws.Range('H8').Value = 'Ensure synthetic_fail contains the column "ENCOUNTER_ID"'
$ws.Range('A9').Value = '847381c1-561e-4a4f-afa7-b2851202b04b'
$ws.Range('B9').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D9').Value = 'Required column SURVEY is missing in synthetic_fail.'
$ws.Range('H9').Value = 'Ensure synthetic_fail contains the column "SURVEY"'
$ws.Range('A10').Value = '1b9a9bef-7f28-43e1-8956-454bc775cf77'
$ws.Range('B10').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D10').Value = 'Required column SURVEY_ID is missing in synthetic_fail.'
$ws.Range('H10').Value = 'Ensure synthetic_fail contains the column "SURVEY_ID"'
$ws.Range('A11').Value = '23df3b39-347d-4091-a2ef-77c6ee3256f8'
$ws.Range('B11').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D11').Value = 'Required column RECORDED_TIME is missing in synthetic_fail.'
$ws.Range('H11').Value = 'Ensure synthetic_fail contains the column "RECORDED_TIME"'
$ws.Range('A12').Value = 'bb2504c7-0bb5-4664-bc49-1a233a0aa321'
$ws.Range('B12').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D12').Value = 'Required column QUESTION is missing in synthetic_fail.'
$ws.Range('H12').Value = 'Ensure synthetic_fail contains the column "QUESTION"'
$ws.Range('A13').Value = '53171127-f04a-420a-9c19-9dc87d3dad33'
$ws.Range('B13').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D13').Value = 'Required column MEAS_VALUE is missing in synthetic_fail.'
$ws.Range('H13').Value = 'Ensure synthetic_fail contains the column "MEAS_VALUE"'
$ws.Range('A14').Value = '42fee73a-c3f5-450e-af0b-84f3a945e3c9'
$ws.Range('B14').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D14').Value = 'Required column QUESTION_CODE is missing in synthetic_fail.'
$ws.Range('H14').Value = 'Ensure synthetic_fail contains the column "QUESTION_CODE"'
$ws.Range('A15').Value = 'db2f4386-3e01-4852-94da-bfcabe69edf9'
$ws.Range('B15').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D15').Value = 'Required column QUESTION_CODE_SYSTEM_NAME is missing in synthetic_fail.'
$ws.Range('H15').Value = 'Ensure synthetic_fail contains the column "QUESTION_CODE_SYSTEM_NAME"'
$ws.Range('A16').Value = '623332c5-6c98-48e6-a314-fff385e13377'
$ws.Range('B16').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D16').Value = 'Required column ANSWER_CODE is missing in synthetic_fail.'
$ws.Range('H16').Value = 'Ensure synthetic_fail contains the column "ANSWER_CODE"'
$ws.Range('A17').Value = '3cc08034-c98a-4bf4-a43d-715c19513897'
$ws.Range('B17').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D17').Value = 'Required column ANSWER_CODE_SYSTEM_NAME is missing in synthetic_fail.'
$ws.Range('H17').Value = 'Ensure synthetic_fail contains the column "ANSWER_CODE_SYSTEM_NAME"'
$ws.Range('A18').Value = '0eb487b4-92a3-4aca-b2db-dd9925dac12b'
$ws.Range('B18').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D18').Value = 'Required column SDOH_DOMAIN is missing in synthetic_fail.'
$ws.Range('H18').Value = 'Ensure synthetic_fail contains the column "SDOH_DOMAIN"'
$ws.Range('A19').Value = '37ad9a18-3b78-4efe-93ca-583262a46cfd'
$ws.Range('B19').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D19').Value = 'Required column NEED_INDICATED is missing in synthetic_fail.'
$ws.Range('H19').Value = 'Ensure synthetic_fail contains the column "NEED_INDICATED"'
$ws.Range('A20').Value = '68eeaf59-3b60-4bf4-89e1-7ea93c100d82'
$ws.Range('B20').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D20').Value = 'Required column VISIT_PART_2_FLAG is missing in synthetic_fail.'
$ws.Range('H20').Value = 'Ensure synthetic_fail contains the column "VISIT_PART_2_FLAG"'
$ws.Range('A21').Value = '37e84512-2d22-44b4-ab8c-140f347f606e'
$ws.Range('B21').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D21').Value = 'Required column VISIT_OMH_FLAG is missing in synthetic_fail.'
$ws.Range('H21').Value = 'Ensure synthetic_fail contains the column "VISIT_OMH_FLAG"'
$ws.Range('A22').Value = 'c7111633-b226-4b9c-b4e6-d7cfc60c527f'
$ws.Range('B22').Value = '6aebd5b8-c1b6-4afb-b160-6328bd174c01'
$ws.Range('D22').Value = 'Required column VISIT_OPWDD_FLAG is missing in synthetic_fail.'
$ws.Range('H22').Value = 'Ensure synthetic_fail contains the column "VISIT_OPWDD_FLAG"'
$ws.Range('A23').Value = 'a37edbfa-bbf2-4a73-8118-653367661f25'
$ws.Range('B23').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A24').Value = 'ae575048-b0e8-4291-8624-46ba902f789a'
$ws.Range('B24').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A25').Value = '77a53eda-afa2-4397-91da-383519fbc6d1'
$ws.Range('B25').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A26').Value = 'd1ac8450-64a2-4d30-9f3f-7fd2aacb844c'
$ws.Range('B26').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A27').Value = 'ce6f461e-3541-477b-a164-076e9a297ef4'
$ws.Range('B27').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A28').Value = 'c8f72ba0-1375-4520-8b44-4150b97064c4'
$ws.Range('B28').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A29').Value = '2b30fb94-3126-47b0-9fe1-d44c10cdf9d1'
$ws.Range('B29').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A30').Value = '6bf5c7b1-f1a3-42a2-98b6-cb77e73a1937'
$ws.Range('B30').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A31').Value = 'e9b098a6-0919-4a58-a56c-5dc3499ec3d8'
$ws.Range('B31').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A32').Value = '3c088bdc-4d7f-44e2-8215-bcffd4f276c3'
$ws.Range('B32').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A33').Value = '26f8f051-6a39-4443-ba4f-a0bd190d4716'
$ws.Range('B33').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A34').Value = '360da663-a72b-4a41-ba76-8a04f4f91087'
$ws.Range('B34').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A35').Value = '52384fea-97e9-4972-b19f-c0be0bbbcddf'
$ws.Range('B35').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A36').Value = 'cac8fdb1-af13-4e47-b49c-341538522788'
$ws.Range('B36').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A37').Value = '25e9ef0e-1757-4ceb-b731-bfe8f6c45f82'
$ws.Range('B37').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A38').Value = '8e3cba8c-8a00-4b58-8cab-9e2a779ea4c1'
$ws.Range('B38').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A39').Value = '58044080-ee49-496f-9250-2f0b7a286b68'
$ws.Range('B39').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A40').Value = '240e19b6-c5e8-46d5-aac8-85eb9e0d9882'
$ws.Range('B40').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A41').Value = 'aac362d8-023c-4b19-9133-753c504b5ae7'
$ws.Range('B41').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A42').Value = 'c426f50f-7047-4f57-b304-8aaec3d8f36f'
$ws.Range('B42').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A43').Value = '973489ba-cfad-4852-b657-4f85963ec05b'
$ws.Range('B43').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A44').Value = '9427ced5-0090-4d82-961f-f5608593becd'
$ws.Range('B44').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A45').Value = '19bd9ff7-1b18-4c9a-b80d-758e90a7a4be'
$ws.Range('B45').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A46').Value = '0285beed-bcbe-442d-8022-417c4c938b58'
$ws.Range('B46').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A47').Value = '53ea0150-5235-4659-897d-9937f29bd2bb'
$ws.Range('B47').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A48').Value = '7ddef889-dd34-4913-9690-5baf7eeef4ee'
$ws.Range('B48').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
$ws.Range('A49').Value = '2d502869-a85e-4f26-928b-12ec67e5d498'
$ws.Range('B49').Value = '9b5ab88a-e757-4520-b89f-b64c440235e1'
